$d = $word.ActiveDocument

# The document has one section whose header/footer slots are split into
# "default" (index 1) and "first page" (index 2) variants:
#   - Footers(1) -> footer2.xml (default footer, Pearson logo, docPr id="2")
#   - Footers(2) -> footer1.xml (first-page footer, Pearson logo, docPr id="3")
#   - Headers(2) -> header1.xml (first-page header, BTEC logo)
# Each holds a single inline picture whose shape name needs to change.

# footer2.xml: Pearson logo, image2.png -> image1.png
$footerDefault = $d.Sections(1).Footers(1)
if ($footerDefault.Exists -and $footerDefault.Range.InlineShapes.Count -gt 0) {
    $footerDefault.Range.InlineShapes(1).Name = "image1.png"
}

# footer1.xml: Pearson logo, image2.png -> image1.png (re-fetch fresh handles
# after the previous mutation so we don't hit a stale-object error)
$footerFirst = $d.Sections(1).Footers(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -gt 0) {
    $footerFirst.Range.InlineShapes(1).Name = "image1.png"
}

# header1.xml: BTEC logo, image1.jpg -> image2.jpg
$headerFirst = $d.Sections(1).Headers(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -gt 0) {
    $headerFirst.Range.InlineShapes(1).Name = "image2.jpg"
}
